$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column F (vGS_min(V)), shifting existing
# columns F..M to G..N.
$ws.Columns.Item(6).Insert()

# New header for the inserted column.
$ws.Cells.Item(1, 6).Value = "Aperture"

# Fill the new column's data rows (2-16) with the value 1, matching the
# number-format style used by column E (scientific notation).
$rng = $ws.Range("F2:F16")
$rng.Value = 1
$rng.NumberFormat = $ws.Range("E2").NumberFormat

# Update the selection to match the post-edit state.
$ws.Range("F2:F16").Select()
